# Clean up "Project 1" / "Project 2" labels on the Electrical sheet so they
# read "Project1" / "Project2" (no space), and leave the Electrical sheet as
# the active/selected sheet (with A3 selected) instead of the Projects sheet.

$wb = $excel.ActiveWorkbook

$wsProjects   = $wb.Worksheets.Item("Projects")
$wsElectrical = $wb.Worksheets.Item("Electrical")

# Remove the stray space in the two project names.
$wsElectrical.Range("A2").Value = "Project1"
$wsElectrical.Range("A3").Value = "Project2"

# Make "Electrical" the active sheet/tab, with A3 selected, matching the
# saved view state of the workbook.
$wsElectrical.Activate()
$wsElectrical.Range("A3").Select()
